# Update the simulated "Return_with_prediction" (G), "return_pct_change" (H)
# and "mean_return_pct_change" (I) columns on Sheet1 with the refreshed
# values produced by the latest run of the auto-recurrence model
# (rows 2-57; column I only carries a value on row 2, the running mean).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "G2" = 0.05284525631479019; "H2" = 58.71840321089577; "I2" = 14.02397192591924
    "G3" = 0.06635313651761979; "H3" = 31.48426382557981
    "G4" = 0.00874539400688538; "H4" = -82.03339224518071
    "G5" = 0.05840167492892269; "H5" = -17.34194212515759
    "G6" = -0.1269181465527413; "H6" = -7.900033141574879
    "G7" = -0.09526340321766807; "H7" = 23.80967834425097
    "G8" = -0.2051602722704262; "H8" = -3.0031062566858
    "G9" = -0.2492974569804259; "H9" = 17.47218421697498
    "G10" = 0.004801985107344942; "H10" = 277.3024573055829
    "G11" = 0.06987018697339868; "H11" = 439.7796066085759
    "G12" = 0.2122349857374143; "H12" = 0.1527564111409068
    "G13" = 0.2495890173009801; "H13" = 6.209633912265371
    "G14" = -0.07548491132566248; "H14" = 17.09534537447127
    "G15" = -0.05531349243235873; "H15" = 22.05934816467435
    "G16" = 0.1721277071627426; "H16" = -10.09414508764179
    "G17" = 0.217817598884789; "H17" = 25.27639661998217
    "G18" = 0.05788591202889228; "H18" = 6.741776625079726
    "G19" = 0.05900508851305674; "H19" = -31.42101122594057
    "G20" = -0.0175531673119257; "H20" = -237.9000281823428
    "G21" = -0.0241411910577574; "H21" = 55.18935767324896
    "G22" = 0.03839384992683196; "H22" = -41.17935007298578
    "G23" = 0.06917175364252126; "H23" = 19.93703650731316
    "G24" = 0.02783624341552676; "H24" = -14.07004798900178
    "G25" = 0.02264039875484063; "H25" = -23.07590903658948
    "G26" = 0.1238873957399087; "H26" = 9.352854754749274
    "G27" = 0.0940073442651724; "H27" = 4.235265564311998
    "G28" = 0.0967929427008974; "H28" = -17.59996730020007
    "G29" = 0.140514323725405; "H29" = 17.45972330193137
    "G30" = 0.06041597327540683; "H30" = -10.1331684573656
    "G31" = 0.06462625349874659; "H31" = -5.816882194061716
    "G32" = 0.04827771598966876; "H32" = 10.55993523788844
    "G33" = 0.1074735038958095; "H33" = 97.78551401662193
    "G34" = -0.007044728849131892; "H34" = 63.11393473692597
    "G35" = 0.05221473187266831; "H35" = 274.2275653586879
    "G36" = -0.007416419952538732; "H36" = -147.9750278627539
    "G37" = 0.03319453452664335; "H37" = 165.0543382154814
    "G38" = 0.06076841127601129; "H38" = -15.29684996146651
    "G39" = 0.02154845398266372; "H39" = -49.96005571711559
    "G40" = 0.05334505902001968; "H40" = 19.26075897279657
    "G41" = 0.02488554277904729; "H41" = 101.3503839117813
    "G42" = 0.089238308656619; "H42" = 70.69635150691185
    "G43" = 0.07048724515789156; "H43" = 41.25732905118962
    "G44" = 0.122626020667023; "H44" = -6.925833284628278
    "G45" = 0.1524901168519333; "H45" = -15.01257698684358
    "G46" = -0.02727374462774964; "H46" = 37.9152226158656
    "G47" = -0.01016103504948941; "H47" = -287.8899017638441
    "G48" = 0.005392831972776166; "H48" = -62.78706406721577
    "G49" = -0.007504836087789579; "H49" = -35.00775699381337
    "G50" = 0.131100575260493; "H50" = -8.277320114208946
    "G51" = 0.1469873229440798; "H51" = 12.23386561341208
    "G52" = 0.0831981119560018; "H52" = 34.29294633119353
    "G53" = 0.05132801040415862; "H53" = -16.10551770674429
    "G54" = -0.1330672448302521; "H54" = -48.99989301210634
    "G55" = -0.08597899886806581; "H55" = 17.1579084097544
    "G56" = 0.1290879263571622; "H56" = -16.78383830258225
    "G57" = 0.1684307244400088; "H57" = 20.75090453998954
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}

Write-Host "Updated $($values.Count) cells across columns G, H and I"
